# Add "cfop" column to the "PI hours" sheet, and a new "cfop hours" sheet
# summarising hours/percentage by cfop value (same pattern as the existing
# "department hours" / "unit(accumulative) hours" sheets).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "PI hours" sheet: new column G = "cfop"
# ---------------------------------------------------------------------
$pi = $wb.Worksheets.Item("PI hours")

# Header (copy formatting from an existing header cell so the new header
# matches: bold, centered, top-aligned, bordered)
$pi.Range("F1").Copy() | Out-Null
$pi.Range("G1").PasteSpecial(-4122) | Out-Null
$pi.Cells.Item(1, 7).Value = "cfop"

$cfopByPI = @(
    "['cfop_HUTCHINSON']",
    "['cfop_KWIAT']",
    "['cfop_NH']",
    "['cfop_GC']",
    "['cfop_MITRA']",
    "['cfop_WORK']",
    "['cfop_CHOUDHURY', 'cfop_RRC']"
)

for ($i = 0; $i -lt $cfopByPI.Length; $i++) {
    $row = $i + 2
    $pi.Cells.Item($row, 7).Value = $cfopByPI[$i]
}

$pi.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. New "cfop hours" sheet, placed after the last existing sheet
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cfopSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$cfopSheet.Name = "cfop hours"

# Header row (same look as the other summary sheets' header rows)
$cfopSheet.Cells.Item(1, 2).Value = "cfop"
$cfopSheet.Cells.Item(1, 3).Value = "hours"
$cfopSheet.Cells.Item(1, 4).Value = "percentage"

# Copy the header formatting (bold, centered, top-aligned, bordered) from
# the "PI hours" sheet's header row rather than rebuilding it property by
# property, so we reuse the existing style instead of minting a new one.
$pi.Range("B1:D1").Copy() | Out-Null
$cfopSheet.Range("B1:D1").PasteSpecial(-4122) | Out-Null
$pi.Application.CutCopyMode = $false

$cfopData = @(
    @{ Name = "cfop_HUTCHINSON"; Hours = 106;  Pct = 44.258872651357 },
    @{ Name = "cfop_GC";         Hours = 51;   Pct = 21.29436325678497 },
    @{ Name = "cfop_KWIAT";      Hours = 41;   Pct = 17.11899791231733 },
    @{ Name = "cfop_NH";         Hours = 35;   Pct = 14.61377870563674 },
    @{ Name = "cfop_MITRA";      Hours = 4;    Pct = 1.670146137787056 },
    @{ Name = "cfop_RRC";        Hours = 1;    Pct = 0.4175365344467641 },
    @{ Name = "cfop_WORK";       Hours = 1;    Pct = 0.4175365344467641 },
    @{ Name = "cfop_CHOUDHURY";  Hours = 0.5;  Pct = 0.208768267223382 }
)

for ($i = 0; $i -lt $cfopData.Length; $i++) {
    $row = $i + 2
    $cfopSheet.Cells.Item($row, 1).Value = $i
    $cfopSheet.Cells.Item($row, 2).Value = $cfopData[$i].Name
    $cfopSheet.Cells.Item($row, 3).Value = $cfopData[$i].Hours
    $cfopSheet.Cells.Item($row, 4).Value = $cfopData[$i].Pct
}

# Column A on the summary sheets carries the same bold/border/centered
# style as the header - copy it from the existing "PI hours" sheet too.
$pi.Range("A2").Copy() | Out-Null
$cfopSheet.Range("A2:A9").PasteSpecial(-4122) | Out-Null
$pi.Application.CutCopyMode = $false

$cfopSheet.PageSetup.LeftMargin = $pi.PageSetup.LeftMargin
$cfopSheet.PageSetup.RightMargin = $pi.PageSetup.RightMargin
$cfopSheet.PageSetup.TopMargin = $pi.PageSetup.TopMargin
$cfopSheet.PageSetup.BottomMargin = $pi.PageSetup.BottomMargin
$cfopSheet.PageSetup.HeaderMargin = $pi.PageSetup.HeaderMargin
$cfopSheet.PageSetup.FooterMargin = $pi.PageSetup.FooterMargin
